# Apply the "speakers" -> poetry content refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab (speakers -> Speakers) ---------------------
$ws.Name = "Speakers"

# --- Replace the header row ------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Image"
$ws.Range("C1").Value = "Url"
$ws.Range("D1").Value = "Title"

# --- Row 2: Emily Dickinson -------------------------------------------------
$ws.Range("A2").Value = "Emily Dickinson"
$ws.Range("B2").Value = "ed.png"
$ws.Range("C2").Value = "https://en.wikipedia.org/wiki/Emily_Dickinson"
$ws.Range("D2").Value = "Success is counted sweetest"

# --- Row 3: Gerard Manley Hopkins ------------------------------------------
$ws.Range("A3").Value = "Gerard Manley Hopkins"
$ws.Range("B3").Value = "gmh.jpg"
$ws.Range("C3").Value = "https://en.wikipedia.org/wiki/Gerard_Manley_Hopkins"
$ws.Range("D3").Value = "The Windhover"

# --- Formatting tweak: C2 picks up its own cell style (new cellXfs entry
#     with an explicit fill applied) -----------------------------------
$ws.Range("C2").Interior.ColorIndex = -4142

# --- Update the saved selection / active cell -------------------------------
$ws.Range("H12").Select()
